$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 5 new mailing-list entries below the existing data (rows 97-101) ---
$ws.Range("A97").Value = "eglantine.hector@curie.fr"
$ws.Range("A98").Value = "joseph.josephides@curie.fr"
$ws.Range("A99").Value = "arnaud.meng@curie.fr"
$ws.Range("A100").Value = "Marion.Salou@curie.fr"
$ws.Range("A101").Value = "nanour.sirab@curie.fr"

# A99 is a "new contributor" marker row, styled like the earlier marker rows
# (e.g. A96) - copy that formatting (font/color) across and match its row height.
$null = $ws.Range("A96").Copy()
$null = $ws.Range("A99").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows(99).RowHeight = 21

# --- Update the view so the newly-added last row is visible/selected ---
$null = $ws.Range("A101").Select()
